# Update MockHeaders test data (formerly ScopeProperties test data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Project/Class names change; Claims row -> Headers row
$ws.Range("A2").Value = "MockHeadersApi"
$ws.Range("B2").Value = "MockHeadersController"
$ws.Range("E2").Value = "A"
$ws.Range("F2").Value = "Headers"
$ws.Range("G2").Value = "header*hdr1=ABC&header*hdr2=DEF"

# Row 3: Expected JSON for scenario A
$ws.Range("A3").Value = "MockHeadersApi"
$ws.Range("B3").Value = "MockHeadersController"
$ws.Range("E3").Value = "A"
$ws.Range("F3").Value = "Expected"
$ws.Range("G3").Value = '[{"Key":"Host","Value":"localhost"},{"Key":"hdr1","Value":"ABC"},{"Key":"hdr2","Value":"DEF"},{"Key":"X-User","Value":"moe@stooges.org"},{"Key":"X-Role","Value":"admin"},{"Key":"X-Role","Value":"user"}]'

# Row 4: Headers for scenario B
$ws.Range("A4").Value = "MockHeadersApi"
$ws.Range("B4").Value = "MockHeadersController"
$ws.Range("E4").Value = "B"
$ws.Range("F4").Value = "Headers"
$ws.Range("G4").Value = "header*X-User=jill&header*X-Role=user"

# Row 5: Expected JSON for scenario B
$ws.Range("A5").Value = "MockHeadersApi"
$ws.Range("B5").Value = "MockHeadersController"
$ws.Range("E5").Value = "B"
$ws.Range("F5").Value = "Expected"
$ws.Range("G5").Value = '[{"Key":"Host","Value":"localhost"},{"Key":"X-User","Value":"larry@stooges.org"},{"Key":"X-Role","Value":"admin"},{"Key":"X-Role","Value":"user"}]'

# Row 6: Headers for scenario C
$ws.Range("A6").Value = "MockHeadersApi"
$ws.Range("B6").Value = "MockHeadersController"
$ws.Range("E6").Value = "C"
$ws.Range("F6").Value = "Headers"
$ws.Range("G6").Value = "header*X-User=jill&header*X-Role=user"

# Row 7: Expected JSON for scenario C
$ws.Range("A7").Value = "MockHeadersApi"
$ws.Range("B7").Value = "MockHeadersController"
$ws.Range("E7").Value = "C"
$ws.Range("F7").Value = "Expected"
$ws.Range("G7").Value = '[{"Key":"Host","Value":"localhost"},{"Key":"X-User","Value":"curly@stooges.org"},{"Key":"X-Role","Value":"readonly"}]'

# Rows 8-10 (old scenario C Claims/Headers/Expected) are no longer needed
$ws.Rows("8:10").Delete()

# Update selection to match the saved view state
$ws.Range("G8").Select()
